# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "68.640.44"
Set-TextCell "E2" "  +1.09%  "
Set-TextCell "D3" "3.278.67"
Set-TextCell "E3" "  +0.35%  "
Set-TextCell "E4" "  +0.02%  "
Set-TextCell "D5" "583.49"
Set-TextCell "E5" "  +0.37%  "
Set-TextCell "D6" "185.55"
Set-TextCell "E6" "  +1.89%  "
Set-TextCell "E7" "  +0.02%  "
Set-TextCell "E8" "  -0.69%  "
Set-TextCell "E9" "  -0.20%  "
Set-TextCell "D10" "6.67"
Set-TextCell "E10" "  -0.96%  "
Set-TextCell "D11" "0.421"
Set-TextCell "E11" "  +0.98%  "
Set-TextCell "D12" "3.854.44"
Set-TextCell "E12" "  +0.53%  "
Set-TextCell "D14" "28.35"
Set-TextCell "E14" "  -0.51%  "
Set-TextCell "D15" "68.684.13"
Set-TextCell "E15" "  +1.19%  "
Set-TextCell "E16" "  +1.40%  "
Set-TextCell "D17" "3.306.15"
Set-TextCell "E17" "  +1.47%  "
Set-TextCell "D18" "5.88"
Set-TextCell "E18" "  +0.64%  "
Set-TextCell "E19" "  +0.91%  "
Set-TextCell "D20" "397.68"
Set-TextCell "E20" "  +5.83%  "
Set-TextCell "D21" "7.76"
Set-TextCell "E21" "  +1.27%  "
Set-TextCell "D22" "71.78"
Set-TextCell "E22" "  +0.88%  "
Set-TextCell "E23" "  -0.01%  "
Set-TextCell "D24" "0.520"
Set-TextCell "E24" "  +1.26%  "
Set-TextCell "E25" "  +0.80%  "
Set-TextCell "E26" "  +4.49%  "
Set-TextCell "D27" "9.78"
Set-TextCell "E27" "  +1.56%  "
Set-TextCell "E28" "  +0.14%  "
Set-TextCell "E29" "  +0.09%  "
Set-TextCell "E30" "  +0.66%  "
Set-TextCell "D31" "23.14"
Set-TextCell "E31" "  +1.29%  "
Set-TextCell "B32" "Aptos"
Set-TextCell "C32" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D32" "7.18"
Set-TextCell "E32" "  +3.59%  "
Set-TextCell "B33" "Fetch.AI"
Set-TextCell "C33" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D33" "1.30"
Set-TextCell "E33" "  +1.96%  "
Set-TextCell "D34" "0.999"
Set-TextCell "E34" "  +0.05%  "
Set-TextCell "D35" "1.52"
Set-TextCell "E35" "  -0.80%  "
Set-TextCell "D36" "163.34"
Set-TextCell "E36" "  +0.44%  "
Set-TextCell "D37" "2.00"
Set-TextCell "E37" "  +7.90%  "
Set-TextCell "D38" "0.827"
Set-TextCell "E38" "  -2.87%  "
Set-TextCell "D39" "26.87"
Set-TextCell "E39" "  -0.24%  "
Set-TextCell "D40" "4.61"
Set-TextCell "E40" "  -0.57%  "
Set-TextCell "D41" "6.62"
Set-TextCell "E41" "  -2.63%  "
Set-TextCell "D42" "2.54"
Set-TextCell "E42" "  -3.05%  "
Set-TextCell "D43" "41.53"
Set-TextCell "E43" "  +1.62%  "
Set-TextCell "D44" "0.0692"
Set-TextCell "E44" "  +1.56%  "
Set-TextCell "D45" "25.53"
Set-TextCell "E45" "  -0.47%  "
Set-TextCell "D46" "2.654.00"
Set-TextCell "E46" "  -0.92%  "
Set-TextCell "D47" "344.25"
Set-TextCell "E47" "  -2.17%  "
Set-TextCell "D48" "0.0283"
Set-TextCell "E48" "  +0.68%  "
Set-TextCell "E49" "  +3.54%  "
Set-TextCell "D50" "31.73"
Set-TextCell "E50" "  +1.85%  "
Set-TextCell "D51" "0.995"
Set-TextCell "E51" "  -0.62%  "
